$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 data (9th fish record)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "HB1603Stn53-6B3I-F1-O1.tif"
$ws.Range("C10").Value = 19.2595
$ws.Range("D10").Value = 9.7287
$ws.Range("E10").Value = 11.2713
$ws.Range("F10").Value = 13.2883
$ws.Range("G10").Value = 15.8586
$ws.Range("AB10").Value = 4

# Update the selection to match the post-edit state
$ws.Range("A10:AB10").Select()
